$d = $word.ActiveDocument

# Locate the paragraph that ends "Working code for each of the three options"
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*of the three options*") {
        $targetPara = $p
    }
}

# Insert a new paragraph after it and fill it with "Responsibilities"
$targetPara.Range.InsertParagraphAfter()
$respPara = $targetPara.Next()
$respPara.Range.Text = "Responsibilities"

# Insert another new paragraph after that one with the Scrum backlog note
$respPara.Range.InsertParagraphAfter()
$backlogPara = $respPara.Next()
$backlogPara.Range.Text = "Look at the Scrum backlog file for this week"

# Make the "Responsibilities" text bold, but not the trailing paragraph mark
# (so the pilcrow/pPr stays un-bolded and boldness doesn't leak onward).
$textRange = $d.Range($respPara.Range.Start, $respPara.Range.End - 1)
$textRange.Font.Bold = 1

# Move the auto "_GoBack" bookmark (tracks last edit point) onto the "R" that
# was just typed - this also removes it from its old location automatically.
$start = $respPara.Range.Start
$bmRange = $d.Range($start, $start + 1)
$d.Bookmarks.Add("_GoBack", $bmRange)
